$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Phase 1 Pre CPP")
$ws2 = $wb.Worksheets.Item("Phase 2")

# New rows of weekly RIHT data to append to "Phase 1 Pre CPP" (rows 2-14).
$data = @(
    @(30414, 0.03507671232876713,  0.09775300847218205,  2061400, 261.6055012161822, 0.01214972872769548),
    @(30421, 0.03515342465753425,  0.02691860832419479,  2061400, 261.7032542246544, 0.01258964950435209),
    @(30428, 0.03523013698630138,  0.01348364722775841,  2061400, 261.7301728329786, 0.01271079211945501),
    @(30435, 0.0353068493150685,   0.01006806763882651,  2061400, 261.7436564802064, 0.01277147297883006),
    @(30442, 0.03538356164383562,  0.005773520527213805, 2061400, 261.7537245478452, 0.01281678260361123),
    @(30477, 0.05007671232876713,  0.5486457479805722,   2061400, 261.7772524173773, 0.01292266577589849),
    @(30484, 0.05015342465753425,  0.1177744491453723,   2061400, 262.3258981653579, 0.01539175258287509),
    @(30491, 0.05023013698630138,  0.03045935413138068,  2061400, 262.4436726145033, 0.01592177644554819),
    @(30498, 0.0503068493150685,   0.01318822766370431,  2061400, 262.4741319686347, 0.01605885358450647),
    @(30505, 0.05038356164383562,  0.01121036431754874,  2061400, 262.4873201962984, 0.01611820495841521),
    @(30512, 0.05046027397260275,  0.01019256660379142,  2061400, 262.4985305606159, 0.01616865529497527),
    @(30519, 0.05053698630136987,  0.01017217838966644,  2061400, 262.5087231272197, 0.01621452520615798),
    @(30526, 0.050613698630137,    0,                     2061400, 262.5188953056094, 0.01626030336365317)
)

$startRow = 2
$endRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
}

# Match the date-column formatting (style used on the "Date" column elsewhere
# in this workbook) by copying it from the existing date cell on "Phase 2".
$ws2.Range("A2").Copy()
$ws1.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)

# The lone data row on "Phase 2" moves to "Phase 1 Pre CPP", so remove it there.
$ws2.Rows.Item(2).Delete()
